$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from J1 to K1, then set header text
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("K1").Value = "kg_quantities_per_portion"

# Populate K2:K254 with kg_quantities_per_portion values
$arr = New-Object 'object[,]' 253,1
$arr[0,0] = 0.01
$arr[1,0] = 0.016666666666666666
$arr[2,0] = 0.03333333333333333
$arr[3,0] = 0.00020537166666666665
$arr[4,0] = 0.13333333333333333
$arr[5,0] = 0.0008214866666666666
$arr[6,0] = 0.015
$arr[7,0] = 0.06666666666666667
$arr[8,0] = 0.049999999999999996
$arr[9,0] = 0.011666666666666667
$arr[10,0] = 0.008333333333333333
$arr[11,0] = 0
$arr[12,0] = 0.03
$arr[13,0] = 0.025
$arr[14,0] = 0.05
$arr[15,0] = 0.00123223
$arr[16,0] = 0.02
$arr[17,0] = 0.15
$arr[18,0] = 0.1
$arr[19,0] = 0.025
$arr[20,0] = 0.025
$arr[21,0] = 0
$arr[22,0] = 0.049999999999999996
$arr[23,0] = 0.03333333333333333
$arr[24,0] = 0.02
$arr[25,0] = 0.049999999999999996
$arr[26,0] = 0.020833333333333332
$arr[27,0] = 0
$arr[28,0] = 0.05
$arr[29,0] = 0.08
$arr[30,0] = 0.025
$arr[31,0] = 0.15
$arr[32,0] = 0.00123223
$arr[33,0] = 0
$arr[34,0] = 0.00492892
$arr[35,0] = 0
$arr[36,0] = 0.0075
$arr[37,0] = 0.1
$arr[38,0] = 0.14
$arr[39,0] = 0.25
$arr[40,0] = 0.015
$arr[41,0] = 0.15
$arr[42,0] = 0
$arr[43,0] = 0.0225
$arr[44,0] = 0.0225
$arr[45,0] = 0.1
$arr[46,0] = 0.1
$arr[47,0] = 0.0075
$arr[48,0] = 0.0075
$arr[49,0] = 0.1
$arr[50,0] = 0.00123223
$arr[51,0] = 0.05
$arr[52,0] = 0.1
$arr[53,0] = 0.025
$arr[54,0] = 0.025
$arr[55,0] = 0.05
$arr[56,0] = 0.0075
$arr[57,0] = 0
$arr[58,0] = 0.00375
$arr[59,0] = 0.05
$arr[60,0] = 0.0075
$arr[61,0] = 0.025
$arr[62,0] = 0.00246446
$arr[63,0] = 0.025
$arr[64,0] = 0.005625
$arr[65,0] = 0.25
$arr[66,0] = 0.125
$arr[67,0] = 0.08
$arr[68,0] = 0.25
$arr[69,0] = 0.03333333333333333
$arr[70,0] = 0.005
$arr[71,0] = 0.0025
$arr[72,0] = 0.016666666666666666
$arr[73,0] = 0.03333333333333333
$arr[74,0] = 0.072
$arr[75,0] = 0.0016429733333333332
$arr[76,0] = 0.0016429733333333332
$arr[77,0] = 0.0008214866666666666
$arr[78,0] = 0.0008214866666666666
$arr[79,0] = 0
$arr[80,0] = 0.0062499999999999995
$arr[81,0] = 0.16666666666666666
$arr[82,0] = 0.016666666666666666
$arr[83,0] = 0.03333333333333333
$arr[84,0] = 0.03333333333333333
$arr[85,0] = 0.03333333333333333
$arr[86,0] = 0.0075
$arr[87,0] = 0.0025
$arr[88,0] = 0.024999999999999998
$arr[89,0] = 0.03333333333333333
$arr[90,0] = 0.016666666666666666
$arr[91,0] = 0.0625
$arr[92,0] = 0.015
$arr[93,0] = 0.05
$arr[94,0] = 0.05
$arr[95,0] = 0.2
$arr[96,0] = 0.05
$arr[97,0] = 0.1
$arr[98,0] = 0.0225
$arr[99,0] = 0.00375
$arr[100,0] = 0.1
$arr[101,0] = 0.005625
$arr[102,0] = 0.025
$arr[103,0] = 0.05
$arr[104,0] = 0.025
$arr[105,0] = 0.000616115
$arr[106,0] = 0.2
$arr[107,0] = 0.025
$arr[108,0] = 0.125
$arr[109,0] = 0.3
$arr[110,0] = 0.0375
$arr[111,0] = 0
$arr[112,0] = 0
$arr[113,0] = 0.1
$arr[114,0] = 0.00375
$arr[115,0] = 0.001875
$arr[116,0] = 0.025
$arr[117,0] = 0.000616115
$arr[118,0] = 0.025
$arr[119,0] = 0.025
$arr[120,0] = 0.0125
$arr[121,0] = 0.05
$arr[122,0] = 0.05
$arr[123,0] = 0
$arr[124,0] = 0
$arr[125,0] = 0.011666666666666667
$arr[126,0] = 0.011666666666666667
$arr[127,0] = 0.16666666666666666
$arr[128,0] = 0.016666666666666666
$arr[129,0] = 0.013333333333333334
$arr[130,0] = 0.0016429733333333332
$arr[131,0] = 0
$arr[132,0] = 0.09999999999999999
$arr[133,0] = 0.041666666666666664
$arr[134,0] = 0
$arr[135,0] = 0.049999999999999996
$arr[136,0] = 0.19999999999999998
$arr[137,0] = 0.008333333333333333
$arr[138,0] = 0.05
$arr[139,0] = 0.0125
$arr[140,0] = 0.05
$arr[141,0] = 0.05
$arr[142,0] = 0.05
$arr[143,0] = 0
$arr[144,0] = 0.25
$arr[145,0] = 0
$arr[146,0] = 0.05
$arr[147,0] = 0.0875
$arr[148,0] = 0.05
$arr[149,0] = 0.0075
$arr[150,0] = 0.1
$arr[151,0] = 0.1
$arr[152,0] = 0.00375
$arr[153,0] = 0.0375
$arr[154,0] = 0.00123223
$arr[155,0] = 0.05
$arr[156,0] = 0.00123223
$arr[157,0] = 0.0625
$arr[158,0] = 0.08
$arr[159,0] = 0
$arr[160,0] = 0.03333333333333333
$arr[161,0] = 0.04666666666666667
$arr[162,0] = 0.06666666666666667
$arr[163,0] = 0.04666666666666667
$arr[164,0] = 0.016666666666666666
$arr[165,0] = 0
$arr[166,0] = 0.10833333333333334
$arr[167,0] = 0.006666666666666667
$arr[168,0] = 0.006666666666666667
$arr[169,0] = 0.0016429733333333332
$arr[170,0] = 0.024999999999999998
$arr[171,0] = 0.03
$arr[172,0] = 0
$arr[173,0] = 0.049999999999999996
$arr[174,0] = 0.049999999999999996
$arr[175,0] = 0
$arr[176,0] = 0.11666666666666665
$arr[177,0] = 0.0016429733333333332
$arr[178,0] = 0.0016429733333333332
$arr[179,0] = 0.0075
$arr[180,0] = 0
$arr[181,0] = 0.09999999999999999
$arr[182,0] = 0.005
$arr[183,0] = 0.26666666666666666
$arr[184,0] = 0.0032859466666666664
$arr[185,0] = 0.049999999999999996
$arr[186,0] = 0.13333333333333333
$arr[187,0] = 0.008333333333333333
$arr[188,0] = 0
$arr[189,0] = 0.00375
$arr[190,0] = 0.0075
$arr[191,0] = 0.005625
$arr[192,0] = 0.2
$arr[193,0] = 0.025
$arr[194,0] = 0.025
$arr[195,0] = 0.0075
$arr[196,0] = 0.03125
$arr[197,0] = 0.025
$arr[198,0] = 0.075
$arr[199,0] = 0.1
$arr[200,0] = 0.2333333333333333
$arr[201,0] = 0.01
$arr[202,0] = 0.03333333333333333
$arr[203,0] = 0.13333333333333333
$arr[204,0] = 0
$arr[205,0] = 0.13333333333333333
$arr[206,0] = 0
$arr[207,0] = 0.0008214866666666666
$arr[208,0] = 0.09999999999999999
$arr[209,0] = 0.05
$arr[210,0] = 0.1
$arr[211,0] = 0.0225
$arr[212,0] = 0.2
$arr[213,0] = 0.1
$arr[214,0] = 0
$arr[215,0] = 0.015
$arr[216,0] = 0
$arr[217,0] = 0.025
$arr[218,0] = 0
$arr[219,0] = 0.05
$arr[220,0] = 0.00375
$arr[221,0] = 0.0075
$arr[222,0] = 0.1
$arr[223,0] = 0.0075
$arr[224,0] = 0.001875
$arr[225,0] = 0.075
$arr[226,0] = 0.05
$arr[227,0] = 0.01125
$arr[228,0] = 0.075
$arr[229,0] = 0.00375
$arr[230,0] = 0
$arr[231,0] = 0
$arr[232,0] = 0.05
$arr[233,0] = 0.00375
$arr[234,0] = 0.015
$arr[235,0] = 0.0075
$arr[236,0] = 0.075
$arr[237,0] = 0.025
$arr[238,0] = 0.00375
$arr[239,0] = 0.075
$arr[240,0] = 0.0125
$arr[241,0] = 0.00375
$arr[242,0] = 0.1
$arr[243,0] = 0.025
$arr[244,0] = 0.075
$arr[245,0] = 0.075
$arr[246,0] = 0.00375
$arr[247,0] = 0
$arr[248,0] = 0
$arr[249,0] = 0.0075
$arr[250,0] = 0.00375
$arr[251,0] = 0.0125
$arr[252,0] = 0.05
$ws.Range("K2:K254").Value = $arr

